# Increment attendee/view counts (column F) for several event rows across
# the "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) worksheets.
# These mirror the same underlying events, so the same counters were bumped
# in each place they are duplicated.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value  = 1177
$wsExpo.Range("F7").Value  = 3842
$wsExpo.Range("F10").Value = 2402
$wsExpo.Range("F17").Value = 96
$wsExpo.Range("F26").Value = 483
$wsExpo.Range("F29").Value = 73
$wsExpo.Range("F33").Value = 869
$wsExpo.Range("F34").Value = 36
$wsExpo.Range("F36").Value = 929
$wsExpo.Range("F37").Value = 1953
$wsExpo.Range("F39").Value = 516
$wsExpo.Range("F43").Value = 1238
$wsExpo.Range("F46").Value = 411

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 63

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 1177
$wsAll.Range("F6").Value  = 3842
$wsAll.Range("F8").Value  = 2402
$wsAll.Range("F14").Value = 96
$wsAll.Range("F23").Value = 483
$wsAll.Range("F26").Value = 63
$wsAll.Range("F29").Value = 73
$wsAll.Range("F33").Value = 869
$wsAll.Range("F34").Value = 36
$wsAll.Range("F36").Value = 929
$wsAll.Range("F37").Value = 1953
$wsAll.Range("F42").Value = 516
$wsAll.Range("F46").Value = 1238
$wsAll.Range("F48").Value = 411
